$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "海盗杰瑞" (Pirate Jerry) to "航海士杰瑞" (Navigator Jerry)
$ws.Range("A5").Value = "航海士杰瑞"

# Normalize full-width parentheses to half-width in existing value-range labels
$ws.Range("H2").Value = "5(6)"
$ws.Range("G4").Value = "400(430)"
$ws.Range("C5").Value = "1.5(2.8)"
$ws.Range("B7").Value = "25(35)"
$ws.Range("D7").Value = "124(174)"
$ws.Range("E8").Value = "2.5(+2.5/+5)"
$ws.Range("B9").Value = "10(25)"
$ws.Range("E13").Value = "2(+25/15s)"
$ws.Range("E21").Value = "2.5(+2)"
$ws.Range("D22").Value = "99(124)"
$ws.Range("F24").Value = "650(500)"
$ws.Range("D31").Value = "75(-9/s*5)"

# Add new column I: "前三分钟推速" (push speed for the first 3 minutes) = base push speed * 0.8
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("H3").Copy()
$ws.Range("I2:I34").PasteSpecial(-4122)

$ws.Range("I1").Value = "前三分钟推速"
$ws.Range("I2").Value = "4(4.8)"
$ws.Range("I3").Formula = "=H3*0.8"
$ws.Range("I4:I34").Formula = "=H4*0.8"
